$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the scaling formula used to compute the hex width (D22):
# change the multiplier from 4 to 4.6. All dependent cells (D23, D25-D29,
# D32-D35, and their ROUND/FLOOR/error/boolean helper columns) recalculate
# automatically from this single change.
$ws.Range("D22").Formula = "=280*4.6"

# Reflect the author's final active selection on the sheet (D22 instead of H13).
$ws.Range("D22").Select()
